$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsProspect = $wb.Worksheets.Item("Prospect")
$wsApplicant = $wb.Worksheets.Item("Applicant")

# ---------------------------------------------------------------------
# Login sheet: swap/update the two hyperlink e-mail addresses
# (A3 first, A2 last -- matches the original authoring order)
# ---------------------------------------------------------------------
$wsLogin.Hyperlinks.Delete()
$wsLogin.Range("A3").Value = "Automation1671@mail.com"

# ---------------------------------------------------------------------
# Prospect sheet: refresh the sample data row (row 2)
# ---------------------------------------------------------------------
$wsProspect.Range("A2").Value = "PK CT"
$wsProspect.Range("B2").Value = "Spring 2020"
$wsProspect.Range("C2").Value = "PK Award"
$wsProspect.Range("D2").Value = "Business"
$wsProspect.Range("E2").Value = "Biology"
$wsProspect.Range("I2").Value = "Junio"

# ---------------------------------------------------------------------
# Applicant sheet: drop the "Campus" column and refresh sample data
# ---------------------------------------------------------------------
$wsApplicant.Columns.Item(3).Delete()
$wsApplicant.Range("B2").Value = "Test"
$wsApplicant.Range("D2").Value = "Event"

# ---------------------------------------------------------------------
# Login sheet: finish the hyperlink swap
# ---------------------------------------------------------------------
$wsLogin.Range("A2").Value = "testselfservice@mail.com"
$wsLogin.Hyperlinks.Add($wsLogin.Range("A3"), "mailto:Automation140@mail.com")
$wsLogin.Hyperlinks.Add($wsLogin.Range("A2"), "mailto:Automation167@mail.com")

# ---------------------------------------------------------------------
# Restore the selections / active sheet seen in the saved workbook
# ---------------------------------------------------------------------
$wsLogin.Activate()
$wsLogin.Range("A4").Select()

$wsApplicant.Activate()
$wsApplicant.Range("C8").Select()

$wsProspect.Activate()
$wsProspect.Range("J2").Select()
